$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '51.991.32'
$ws.Range("E2").Value = '  +1.06%  '
$ws.Range("D3").Value = '2.822.65'
$ws.Range("E3").Value = '  +2.74%  '
$ws.Range("E4").Value = '  -0.02%  '
$ws.Range("D5").Value = "'355.73"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +7.00%  '
$ws.Range("D6").Value = "'113.77"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -1.63%  '
$ws.Range("E7").Value = '  +2.70%  '
$ws.Range("E8").Value = '  -0.01%  '
$ws.Range("E9").Value = '  +4.80%  '
$ws.Range("D10").Value = "'41.91"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +0.83%  '
$ws.Range("D11").Value = "'0.0852"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -0.23%  '
$ws.Range("D12").Value = "'20.09"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -0.31%  '
$ws.Range("E13").Value = '  +1.38%  '
$ws.Range("E14").Value = '  +1.52%  '
$ws.Range("D15").Value = '3.243.63'
$ws.Range("E15").Value = '  +2.16%  '
$ws.Range("D16").Value = '2.833.01'
$ws.Range("E16").Value = '  +3.64%  '
$ws.Range("D17").Value = "'0.898"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +1.71%  '
$ws.Range("D18").Value = '51.840.06'
$ws.Range("E18").Value = '  +0.77%  '
$ws.Range("E19").Value = '  +8.17%  '
$ws.Range("E20").Value = '  -1.28%  '
$ws.Range("D21").Value = "'13.56"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +1.21%  '
$ws.Range("E22").Value = '  +2.67%  '
$ws.Range("D23").Value = "'270.22"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -2.98%  '
$ws.Range("D24").Value = "'69.78"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +0.56%  '
$ws.Range("D25").Value = "'2.80"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +5.74%  '
$ws.Range("D26").Value = "'26.85"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +0.23%  '
$ws.Range("E27").Value = '  +0.04%  '
$ws.Range("E28").Value = '  +1.47%  '
$ws.Range("E29").Value = '  +1.72%  '
$ws.Range("E30").Value = '  -0.42%  '
$ws.Range("B31").Value = 'VeChain'
$ws.Range("C31").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D31").Value = "'0.0460"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +34.22%  '
$ws.Range("B32").Value = 'OKB'
$ws.Range("C32").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D32").Value = "'50.82"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +1.98%  '
$ws.Range("B33").Value = 'InjectiveProtocol'
$ws.Range("C33").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D33").Value = "'33.91"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -3.10%  '
$ws.Range("D34").Value = "'5.85"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +5.76%  '
$ws.Range("E35").Value = '  +0.77%  '
$ws.Range("E36").Value = '  -0.10%  '
$ws.Range("E37").Value = '  +0.28%  '
$ws.Range("D38").Value = "'4.91"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -1.51%  '
$ws.Range("D39").Value = "'3.21"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -0.09%  '
$ws.Range("E40").Value = '  -3.78%  '
$ws.Range("B41").Value = 'Stacks'
$ws.Range("C41").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D41").Value = "'2.58"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +5.38%  '
$ws.Range("B42").Value = 'Monero'
$ws.Range("C42").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D42").Value = "'128.01"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +0.58%  '
$ws.Range("B43").Value = 'EnergySwap'
$ws.Range("C43").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D43").Value = "'23.52"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +2.39%  '
$ws.Range("E44").Value = '  +1.54%  '
$ws.Range("E45").Value = '  +0.66%  '
$ws.Range("E46").Value = '  +1.00%  '
$ws.Range("D47").Value = '2.075.89'
$ws.Range("E47").Value = '  -0.55%  '
$ws.Range("E48").Value = '  +3.83%  '
$ws.Range("B49").Value = 'THORChain'
$ws.Range("C49").Value = 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
$ws.Range("D49").Value = "'5.70"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +3.28%  '
$ws.Range("B50").Value = 'SEI'
$ws.Range("C50").Value = 'https://coinranking.com/coin/8nxCqs-uj+sei-sei'
$ws.Range("D50").Value = "'0.949"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +9.53%  '
$ws.Range("D51").Value = "'60.59"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +1.35%  '
